$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (137) onto the
# two new rows (138, 139) so that the "Indice" column keeps its bold /
# bordered / centered style and the "data_partida" column keeps its
# date-time number format.
$ws.Range("A137:V137").Copy()
$ws.Range("A138:V139").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 138 : FC Voluntari 0 x 0 Petrolul ---
$ws.Cells.Item(138, 1).Value = 137
$ws.Cells.Item(138, 2).Value = "romania"
$ws.Cells.Item(138, 3).Value = "liga-1"
$ws.Cells.Item(138, 4).Value = "2023-2024"
$ws.Cells.Item(138, 5).Value = 45261.66666666666
$ws.Cells.Item(138, 6).Value = "FC Voluntari"
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = "Petrolul"
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 2.51
$ws.Cells.Item(138, 11).Value = "30/11/2023 08:12"
$ws.Cells.Item(138, 12).Value = 2.93
$ws.Cells.Item(138, 13).Value = "01/12/2023 15:59"
$ws.Cells.Item(138, 14).Value = 2.95
$ws.Cells.Item(138, 15).Value = "30/11/2023 08:12"
$ws.Cells.Item(138, 16).Value = 2.82
$ws.Cells.Item(138, 17).Value = "01/12/2023 15:58"
$ws.Cells.Item(138, 18).Value = 2.97
$ws.Cells.Item(138, 19).Value = "30/11/2023 08:12"
$ws.Cells.Item(138, 20).Value = 2.87
$ws.Cells.Item(138, 21).Value = "01/12/2023 15:59"
$ws.Cells.Item(138, 22).Value = "https://www.betexplorer.com/football/romania/liga-1/voluntari-petrolul/2qPCPF03/"

# --- Row 139 : Poli Iasi 2 x 3 Farul Constanta ---
$ws.Cells.Item(139, 1).Value = 138
$ws.Cells.Item(139, 2).Value = "romania"
$ws.Cells.Item(139, 3).Value = "liga-1"
$ws.Cells.Item(139, 4).Value = "2023-2024"
$ws.Cells.Item(139, 5).Value = 45261.79166666666
$ws.Cells.Item(139, 6).Value = "Poli Iasi"
$ws.Cells.Item(139, 7).Value = 2
$ws.Cells.Item(139, 8).Value = "Farul Constanta"
$ws.Cells.Item(139, 9).Value = 3
$ws.Cells.Item(139, 10).Value = 3.32
$ws.Cells.Item(139, 11).Value = "30/11/2023 08:12"
$ws.Cells.Item(139, 12).Value = 3.45
$ws.Cells.Item(139, 13).Value = "01/12/2023 18:53"
$ws.Cells.Item(139, 14).Value = 3.22
$ws.Cells.Item(139, 15).Value = "30/11/2023 08:12"
$ws.Cells.Item(139, 16).Value = 3.3
$ws.Cells.Item(139, 17).Value = "01/12/2023 18:52"
$ws.Cells.Item(139, 18).Value = 2.16
$ws.Cells.Item(139, 19).Value = "30/11/2023 08:12"
$ws.Cells.Item(139, 20).Value = 2.21
$ws.Cells.Item(139, 21).Value = "01/12/2023 18:53"
$ws.Cells.Item(139, 22).Value = "https://www.betexplorer.com/football/romania/liga-1/poli-iasi-farul-constanta/d6OGOZF9/"
